$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26, shifting existing rows 26-50 down to 27-51.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the Emerald exposure site data.
$ws.Range("A26").Value = "Emerald"
$ws.Range("B26").Value = "Puffing Billy Railway Lakeside Station, Emerald Lake Rd, Emerald"
$ws.Range("C26").Value = "31/12/20 3:00pm - 5:00pm"
$ws.Range("D26").Value = "Case vistied venue"
